$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.456.36"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "'1.647.25"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'298.88"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "'0.3543"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "'50.18"
$ws.Range("D10").Value = "'0.08087"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'22.09"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "'6.398"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "'7.348"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'0.00001201"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").Value = "'1.644.66"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'0.06950"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'6.770"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'17.39"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'12.45"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").Value = "'23.480.47"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "'2.509"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "'2.893"
$ws.Range("E26").Value = "  -6.67%  "
$ws.Range("D27").Value = "'20.88"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "'152.14"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "'5.204"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").Value = "'132.81"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "'1.833.64"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'6.942"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "'2.145"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").Value = "'0.9911"
$ws.Range("E35").Value = "  -8.43%  "
$ws.Range("D36").Value = "'0.02718"
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("D37").Value = "'0.08756"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'0.2440"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").Value = "'5.938"
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'0.06792"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").Value = "'0.6897"
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'15.76"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.296"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'0.6365"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").Value = "'2.255"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "'3.911"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "'0.07726"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").Value = "'127.84"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'1.151"
$ws.Range("E51").Value = "  -3.88%  "
